$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 18:12"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 6267472
$ws.Range("C4").Value = 9901
$ws.Range("D4").Value = 3500024
$ws.Range("E4").Value = 2578256
$ws.Range("G4").Value = 292
$ws.Range("H4").Value = 189192

# Row 6: India -> India
$ws.Range("B6").Value = 3823449
$ws.Range("C6").Value = 57341
$ws.Range("D6").Value = 2946920
$ws.Range("E6").Value = 809374
$ws.Range("G6").Value = 695
$ws.Range("H6").Value = 67155

# Row 16: Reino Unido -> Reino Unido
$ws.Range("B16").Value = 338676
$ws.Range("C16").Value = 1508
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 41514

# Row 22: Italia -> Italia
$ws.Range("B22").Value = 271515
$ws.Range("C22").Value = 1326
$ws.Range("D22").Value = 208201
$ws.Range("E22").Value = 27817
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 35497

# Row 23: Alemania -> Alemania
$ws.Range("B23").Value = 246702
$ws.Range("C23").Value = 701
$ws.Range("E23").Value = 15517

# Row 27: Canada -> Canada
$ws.Range("B27").Value = 129691
$ws.Range("C27").Value = 266
$ws.Range("D27").Value = 114818
$ws.Range("E27").Value = 5739
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 9134

# Row 52: Singapur -> Singapur
$ws.Range("D52").Value = 55891
$ws.Range("E52").Value = 983

# Row 91: Croacia -> Grecia
$ws.Range("A91").Value = "Grecia"
$ws.Range("B91").Value = 10757
$ws.Range("C91").Value = 233
$ws.Range("D91").Value = 3804
$ws.Range("E91").Value = 6680
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 273

# Row 92: Grecia -> Croacia
$ws.Range("A92").Value = "Croacia"
$ws.Range("B92").Value = 10725
$ws.Range("C92").Value = 311
$ws.Range("D92").Value = 7968
$ws.Range("E92").Value = 2566
$ws.Range("G92").Value = 4
$ws.Range("H92").Value = 191

# Row 94: Albania -> Albania
$ws.Range("B94").Value = 9728
$ws.Range("C94").Value = 122
$ws.Range("D94").Value = 5582
$ws.Range("E94").Value = 3850
$ws.Range("G94").Value = 6
$ws.Range("H94").Value = 296

# Row 101: Finlandia -> Finlandia
$ws.Range("D101").Value = 7350
$ws.Range("E101").Value = 475

# Row 103: Namibia -> Namibia
$ws.Range("B103").Value = 7844
$ws.Range("C103").Value = 152
$ws.Range("D103").Value = 3454
$ws.Range("E103").Value = 4308
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 82

# Row 117: Surinam -> Mozambique
$ws.Range("A117").Value = "Mozambique"
$ws.Range("B117").Value = 4117
$ws.Range("C117").Value = 78
$ws.Range("D117").Value = 2170
$ws.Range("E117").Value = 1922
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 25

# Row 118: Cuba -> Surinam
$ws.Range("A118").Value = "Surinam"
$ws.Range("B118").Value = 4089
$ws.Range("D118").Value = 3171
$ws.Range("E118").Value = 846
$ws.Range("H118").Value = 72

# Row 119: Eslovaquia -> Cuba
$ws.Range("A119").Value = "Cuba"
$ws.Range("B119").Value = 4065
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 3395
$ws.Range("E119").Value = 575
$ws.Range("H119").Value = 95

# Row 120: Mozambique -> Eslovaquia
$ws.Range("A120").Value = "Eslovaquia"
$ws.Range("B120").Value = 4042
$ws.Range("C120").Value = 53
$ws.Range("D120").Value = 2523
$ws.Range("E120").Value = 1486
$ws.Range("H120").Value = 33

# Row 128: Sri Lanka -> Sri Lanka
$ws.Range("B128").Value = 3101
$ws.Range("C128").Value = 9
$ws.Range("E128").Value = 206

# Row 129: Gambia -> Gambia
$ws.Range("B129").Value = 3067
$ws.Range("C129").Value = 38
$ws.Range("E129").Value = 1938
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 97

# Row 147: Trinidad yTobago -> Trinidad yTobago
$ws.Range("B147").Value = 1839
$ws.Range("C147").Value = 42
$ws.Range("D147").Value = 691
$ws.Range("E147").Value = 1120
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 28

# Row 153: Republica de Chipre -> Republica de Chipre
$ws.Range("B153").Value = 1495
$ws.Range("C153").Value = 5
$ws.Range("E153").Value = 335

# Row 166: Birmania -> Birmania
$ws.Range("B166").Value = 995
$ws.Range("C166").Value = 76
$ws.Range("D166").Value = 358
$ws.Range("E166").Value = 631
